# The deck ships with two embedded themes:
#   ppt/theme/theme1.xml -> "Integral"      (used by the slide master)
#   ppt/theme/theme2.xml -> "Office Theme"  (used by the notes master)
#
# The authored edit swaps the two themes' contents: the slide master's
# theme becomes the stock "Office Theme" palette, and the notes master's
# theme becomes the old "Integral" palette. The only part of that swap
# reachable through the PowerPoint object model is the slide master's
# theme colour scheme (12 colour slots, in clrScheme document order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) - so drive the
# "Integral" -> "Office Theme" colour change through there.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

function RgbValue($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme palette (RRGGBB), applied in clrScheme slot order.
$colors.Item(1).RGB  = RgbValue 0x00 0x00 0x00   # dk1      000000
$colors.Item(2).RGB  = RgbValue 0xFF 0xFF 0xFF   # lt1      FFFFFF
$colors.Item(3).RGB  = RgbValue 0x44 0x54 0x6A   # dk2      44546A
$colors.Item(4).RGB  = RgbValue 0xE7 0xE6 0xE6   # lt2      E7E6E6
$colors.Item(5).RGB  = RgbValue 0x5B 0x9B 0xD5   # accent1  5B9BD5
$colors.Item(6).RGB  = RgbValue 0xED 0x7D 0x31   # accent2  ED7D31
$colors.Item(7).RGB  = RgbValue 0xA5 0xA5 0xA5   # accent3  A5A5A5
$colors.Item(8).RGB  = RgbValue 0xFF 0xC0 0x00   # accent4  FFC000
$colors.Item(9).RGB  = RgbValue 0x44 0x72 0xC4   # accent5  4472C4
$colors.Item(10).RGB = RgbValue 0x70 0xAD 0x47   # accent6  70AD47
$colors.Item(11).RGB = RgbValue 0x05 0x63 0xC1   # hlink    0563C1
$colors.Item(12).RGB = RgbValue 0x95 0x4F 0x72   # folHlink 954F72
